$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing electric-related header labels ---
$ws.Range("I1").Value = "Electric Choice ID"
$ws.Range("J1").Value = "Electric Rate Code"
$ws.Range("L1").Value = "Electric Usage (kWh)"

# --- Add the new gas header columns (M, N, O), matching the style of the ---
# --- existing header row so formatting (bold, border, centered) carries over ---
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("M1").Value = "Gas Choice ID"
$ws.Range("N1").Value = "Gas Rate Code"
$ws.Range("O1").Value = "Gas Usage (therms)"

# --- Extend row 2 with matching (blank) data cells under the new columns ---
# A literal "" assignment gets dropped by the save pipeline (treated as
# clearing the cell), so force real, empty text cells the same way Excel's
# UI does for a bare leading apostrophe, then strip the resulting
# quote-prefix style so the cell ends up plain/unstyled -- matching B2:L2.
$ws.Range("M2").Value = "'"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'"
$ws.Range("M2:O2").Style = "Normal"
